$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 12211.677
$ws.Range("I15").Value = 12211.677
$ws.Range("K15").Value = 36635.031
$ws.Range("M15").Value = -36466.031
$ws.Range("H17").Value = 3731200.5
$ws.Range("J17").Value = 3731200.5
$ws.Range("L17").Value = 11193601.5
$ws.Range("N17").Value = -11193937.5
$ws.Range("H40").Value = 1835.6212
$ws.Range("I40").Value = 1874.4615
$ws.Range("J40").Value = 1691.3572
$ws.Range("K40").Value = 1874.4615
$ws.Range("L40").Value = 1691.3572
$ws.Range("M40").Value = -1699.4615
$ws.Range("N40").Value = -2041.3572
$ws.Range("H64").Value = 3875
$ws.Range("I64").Value = 4142.857
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 4142.857
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -3894.857
$ws.Range("N64").Value = -3996
$ws.Range("H67").Value = 3875
$ws.Range("I67").Value = 4142.857
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 4142.857
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -3284.857
$ws.Range("N67").Value = -5216
$ws.Range("H112").Value = 31261896
$ws.Range("I112").Value = 685
$ws.Range("J112").Value = 41682300
$ws.Range("K112").Value = 2055
$ws.Range("L112").Value = 125046900
$ws.Range("M112").Value = -947
$ws.Range("N112").Value = -125049116
$ws.Range("H129").Value = 886.67
$ws.Range("J129").Value = 975.8554
$ws.Range("L129").Value = 2927.5662
$ws.Range("N129").Value = -12927.5662
$ws.Range("H132").Value = 9805412
$ws.Range("I132").Value = 1114.3889
$ws.Range("K132").Value = 3343.1667
$ws.Range("M132").Value = -813.1666999999998
$ws.Range("H137").Value = 1344.081
$ws.Range("I137").Value = 1228.2593
$ws.Range("K137").Value = 3684.7779
$ws.Range("M137").Value = -1134.7779

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6015.1865
$ws.Range("I32").Value = 4332.8066
$ws.Range("K32").Value = 4332.8066
$ws.Range("M32").Value = -4045.8066

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 51581.6
$ws.Range("I20").Value = 2004
$ws.Range("J20").Value = 84633.336
$ws.Range("K20").Value = 2004
$ws.Range("L20").Value = 84633.336
$ws.Range("M20").Value = -1757
$ws.Range("N20").Value = -85127.336
$ws.Range("H63").Value = 42200
$ws.Range("J63").Value = 42200
$ws.Range("L63").Value = 42200
$ws.Range("N63").Value = -43572
$ws.Range("H66").Value = 42200
$ws.Range("J66").Value = 42200
$ws.Range("L66").Value = 126600
$ws.Range("N66").Value = -133464
$ws.Range("H99").Value = 58824500
$ws.Range("I99").Value = 66667564
$ws.Range("K99").Value = 66667564
$ws.Range("M99").Value = -66666066
$ws.Range("H107").Value = 1048.3462
$ws.Range("I107").Value = 1111.8572
$ws.Range("J107").Value = 781.6
$ws.Range("K107").Value = 1111.8572
$ws.Range("L107").Value = 781.6
$ws.Range("M107").Value = 808.1428000000001
$ws.Range("N107").Value = -4621.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9618410
$ws.Range("I31").Value = 1641
$ws.Range("J31").Value = 26321218
$ws.Range("K31").Value = 1641
$ws.Range("L31").Value = 26321218
$ws.Range("M31").Value = -1346
$ws.Range("N31").Value = -26321808
$ws.Range("H34").Value = 9618410
$ws.Range("I34").Value = 1641
$ws.Range("J34").Value = 26321218
$ws.Range("K34").Value = 1641
$ws.Range("L34").Value = 26321218
$ws.Range("M34").Value = -1439
$ws.Range("N34").Value = -26321622
$ws.Range("H58").Value = 6411654
$ws.Range("I58").Value = 11905696
$ws.Range("J58").Value = 1938.75
$ws.Range("K58").Value = 11905696
$ws.Range("L58").Value = 1938.75
$ws.Range("M58").Value = -11905493
$ws.Range("N58").Value = -2344.75
$ws.Range("H94").Value = 3286.9644
$ws.Range("I94").Value = 3714.182
$ws.Range("J94").Value = 3010.5293
$ws.Range("K94").Value = 3714.182
$ws.Range("L94").Value = 3010.5293
$ws.Range("M94").Value = -3263.182
$ws.Range("N94").Value = -3912.5293
$ws.Range("H99").Value = 3909884.2
$ws.Range("I99").Value = 2948
$ws.Range("J99").Value = 17863228
$ws.Range("K99").Value = 2948
$ws.Range("L99").Value = 17863228
$ws.Range("M99").Value = -1450
$ws.Range("N99").Value = -17866224
$ws.Range("H105").Value = 20835242
$ws.Range("I105").Value = 25642876
$ws.Range("K105").Value = 25642876
$ws.Range("M105").Value = -25641129
$ws.Range("H126").Value = 3909884.2
$ws.Range("I126").Value = 2948
$ws.Range("J126").Value = 17863228
$ws.Range("K126").Value = 8844
$ws.Range("L126").Value = 53589684
$ws.Range("M126").Value = -6374
$ws.Range("N126").Value = -53594624
$ws.Range("H134").Value = 11495717
$ws.Range("I134").Value = 13890328
$ws.Range("J134").Value = 1582.8
$ws.Range("K134").Value = 41670984
$ws.Range("L134").Value = 4748.4
$ws.Range("M134").Value = -41668449
$ws.Range("N134").Value = -9818.4
$ws.Range("H136").Value = 6411654
$ws.Range("I136").Value = 11905696
$ws.Range("J136").Value = 1938.75
$ws.Range("K136").Value = 35717088
$ws.Range("L136").Value = 5816.25
$ws.Range("M136").Value = -35714538
$ws.Range("N136").Value = -10916.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 5882446.5
$ws.Range("J23").Value = 97
$ws.Range("L23").Value = 291
$ws.Range("N23").Value = -761
$ws.Range("H70").Value = 1943.8572
$ws.Range("I70").Value = 746
$ws.Range("K70").Value = 2238
$ws.Range("M70").Value = -1923
$ws.Range("H73").Value = 1943.8572
$ws.Range("I73").Value = 746
$ws.Range("K73").Value = 2238
$ws.Range("M73").Value = -1146
$ws.Range("H75").Value = 20408552
$ws.Range("I75").Value = 146.2
$ws.Range("K75").Value = 438.6
$ws.Range("M75").Value = 559.4000000000001
$ws.Range("H78").Value = 20408552
$ws.Range("I78").Value = 146.2
$ws.Range("K78").Value = 1315.8
$ws.Range("M78").Value = 3676.2
$ws.Range("H86").Value = 675.7273
$ws.Range("I86").Value = 666.25
$ws.Range("J86").Value = 701
$ws.Range("K86").Value = 1998.75
$ws.Range("L86").Value = 2103
$ws.Range("M86").Value = -812.75
$ws.Range("N86").Value = -4475
$ws.Range("H87").Value = 1379.8
$ws.Range("I87").Value = 1379.8
$ws.Range("K87").Value = 4139.4
$ws.Range("M87").Value = -2891.4
$ws.Range("H89").Value = 675.7273
$ws.Range("I89").Value = 666.25
$ws.Range("J89").Value = 701
$ws.Range("K89").Value = 5996.25
$ws.Range("L89").Value = 6309
$ws.Range("M89").Value = -68.25
$ws.Range("N89").Value = -18165
$ws.Range("H90").Value = 1379.8
$ws.Range("I90").Value = 1379.8
$ws.Range("K90").Value = 12418.2
$ws.Range("M90").Value = -6178.199999999999
$ws.Range("H92").Value = 666.6667
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 850
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 2550
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -5046
$ws.Range("H104").Value = 3500
$ws.Range("J104").Value = 3500
$ws.Range("L104").Value = 10500
$ws.Range("N104").Value = -15742
$ws.Range("H113").Value = 2400492.8
$ws.Range("I113").Value = 2632070.5
$ws.Range("J113").Value = 1667163.4
$ws.Range("K113").Value = 7896211.5
$ws.Range("L113").Value = 5001490.199999999
$ws.Range("M113").Value = -7894041.5
$ws.Range("N113").Value = -5005830.199999999
$ws.Range("H119").Value = 215468.25
$ws.Range("I119").Value = 2764.5
$ws.Range("K119").Value = 8293.5
$ws.Range("M119").Value = -3455.5
$ws.Range("H121").Value = 1002.65
$ws.Range("I121").Value = 676
$ws.Range("J121").Value = 1049.3143
$ws.Range("K121").Value = 2028
$ws.Range("L121").Value = 3147.9429
$ws.Range("M121").Value = -718
$ws.Range("N121").Value = -5767.9429
$ws.Range("H131").Value = 2942065
$ws.Range("I131").Value = 10000416
$ws.Range("J131").Value = 1085.625
$ws.Range("K131").Value = 30001248
$ws.Range("L131").Value = 3256.875
$ws.Range("M131").Value = -29996208
$ws.Range("N131").Value = -13336.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 9868.691999999999
$ws.Range("I126").Value = 12099.2
$ws.Range("J126").Value = 2433.6667
$ws.Range("K126").Value = 36297.60000000001
$ws.Range("L126").Value = 7301.000100000001
$ws.Range("M126").Value = -33827.60000000001
$ws.Range("N126").Value = -12241.0001
$ws.Range("H132").Value = 11908080
$ws.Range("I132").Value = 16669590
$ws.Range("K132").Value = 50008770
$ws.Range("M132").Value = -50006240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11181395
$ws.Range("I132").Value = 14325225
$ws.Range("J132").Value = 3331.6667
$ws.Range("K132").Value = 42975675
$ws.Range("L132").Value = 9995.000100000001
$ws.Range("M132").Value = -42973145
$ws.Range("N132").Value = -15055.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 970.82355
$ws.Range("I113").Value = 741.1
$ws.Range("J113").Value = 1299
$ws.Range("K113").Value = 2223.3
$ws.Range("L113").Value = 3897
$ws.Range("M113").Value = -53.30000000000018
$ws.Range("N113").Value = -8237
$ws.Range("H136").Value = 3334260
$ws.Range("I136").Value = 531.34375
$ws.Range("J136").Value = 9260889
$ws.Range("K136").Value = 1594.03125
$ws.Range("L136").Value = 27782667
$ws.Range("M136").Value = 955.96875
$ws.Range("N136").Value = -27787767
